$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Variables List Block" upper-left/lower-right index cells
# (D5:D11) from the *43 scenario row references to the *46 ones.
$ws.Range("D5").Value  = "A46"
$ws.Range("D6").Value  = "B46"
$ws.Range("D7").Value  = "C46"
$ws.Range("D8").Value  = "G46"
$ws.Range("D9").Value  = "H46"
$ws.Range("D10").Value = "I46"
$ws.Range("D11").Value = "J46"

# Demands List Block lower-right cell (D15) stays "E261" - unchanged value.
$ws.Range("D15").Value = "E261"

# Update the active selection to reflect the edited block.
$ws.Range("D5:D11").Select()
